# feat: add 2022-Q1 data
#
# This adds a new sheet "2022-Q1" (fund-level holdings detail, same shape as
# the "2020-Q4" / "2021-Q1" sheets) positioned right before the "总计"
# (summary) sheet, and updates the "总计" sheet with a new top row summarizing
# the 2022-Q1 quarter, pushing the existing rows down by one.

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("2021-Q1")

# --- 1. Create the new "2022-Q1" sheet by duplicating an existing detail
#        sheet (keeps sheetPr/styles/pageMargins identical) and inserting it
#        right before "总计".
# NOTE: worksheet variables in this runtime resolve by position, not
# identity - once the new sheet is spliced in, a previously-captured
# reference to "总计" would silently start pointing at the new sheet
# instead. So "总计" is (re)fetched by name only *after* the sheet
# list has settled into its final shape.
$template.Copy($wb.Worksheets.Item("总计"))
$newSheet = $wb.Worksheets.Item($template.Index + 1)
$newSheet.Name = "2022-Q1"

$total = $wb.Worksheets.Item("总计")

# --- 2. Fill in the header row for the new sheet.
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# --- 3. Fill in the two fund holding rows. Codes/names/numbers are kept as
#        text (matching the sheet's existing convention), so force text
#        formatting before assigning via NumberFormat "@" and resetting the
#        style back to Normal afterwards so no stray style gets attached.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$newSheet.Range("A2").Value = 0
Set-TextValue $newSheet.Range("B2") "001411"
Set-TextValue $newSheet.Range("C2") "诺安创新驱动灵活配置混合A"
Set-TextValue $newSheet.Range("D2") "3.96"
Set-TextValue $newSheet.Range("E2") "94.02"
Set-TextValue $newSheet.Range("F2") "4.47"
Set-TextValue $newSheet.Range("G2") "0.1770"
$newSheet.Range("H2").Value = 8

# A3 is a brand new cell (the template only had one data row), so first
# copy A2's style (bold/bordered index-column look) down onto it.
$newSheet.Range("A2").Copy($newSheet.Range("A3"))
$newSheet.Range("A3").Value = 1
Set-TextValue $newSheet.Range("B3") "002051"
Set-TextValue $newSheet.Range("C3") "诺安创新驱动灵活配置混合C"
Set-TextValue $newSheet.Range("D3") "1.33"
Set-TextValue $newSheet.Range("E3") "94.02"
Set-TextValue $newSheet.Range("F3") "4.47"
Set-TextValue $newSheet.Range("G3") "0.0595"
$newSheet.Range("H3").Value = 8

# --- 4. Update the "总计" summary sheet: push the existing rows down by one
#        (reading with Value2, since Value has a read quirk in this runtime)
#        and insert the new 2022-Q1 summary row at the top of the data.
$oldRow2B = $total.Range("B2").Value2
$oldRow2C = $total.Range("C2").Value2
$oldRow2D = $total.Range("D2").Value2
$oldRow3B = $total.Range("B3").Value2
$oldRow3C = $total.Range("C3").Value2
$oldRow3D = $total.Range("D3").Value2

# Row 4 (was row 3): copy A3's style down to the freshly used A4 cell.
$total.Range("A3").Copy($total.Range("A4"))
$total.Range("A4").Value = 2
$total.Range("B4").Value = $oldRow3B
$total.Range("C4").Value = $oldRow3C
$total.Range("D4").Value = $oldRow3D

# Row 3 (was row 2).
$total.Range("A3").Value = 1
$total.Range("B3").Value = $oldRow2B
$total.Range("C3").Value = $oldRow2C
$total.Range("D3").Value = $oldRow2D

# Row 2: the brand-new 2022-Q1 summary row.
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.24

# --- 5. Restore the originally active tab (creating/copying sheets along
#        the way shifts which sheet is "active"; the source workbook has
#        "2020-Q4" selected).
$wb.Worksheets.Item("2020-Q4").Activate()
